$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.343.20'
$ws.Range("E2").Value = '  -6.78%  '
$ws.Range("D3").Value = '2.442.86'
$ws.Range("E3").Value = '  -10.84%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '468.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.994'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.490'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.60%  '
$ws.Range("D9").Value = '2.464.91'
$ws.Range("E9").Value = '  -10.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0962'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.35'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.319'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.122'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.50%  '
$ws.Range("D14").Value = '2.865.80'
$ws.Range("E14").Value = '  -11.30%  '
$ws.Range("D15").Value = '54.182.75'
$ws.Range("E15").Value = '  -7.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000134'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.55%  '
$ws.Range("D18").Value = '2.461.62'
$ws.Range("E18").Value = '  -9.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -8.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '314.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -12.04%  '
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -11.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '56.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.45%  '
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.387'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.65%  '
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '2.540.76'
$ws.Range("E28").Value = '  -12.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.154'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("E31").Value = '  -0.31%  '
$ws.Range("D32").Value = '0.0₃0730'
$ws.Range("E32").Value = '  -7.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '150.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.82'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.59'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -12.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.07'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.808'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.994'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.605'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0532'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.72%  '
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.83%  '
$ws.Range("D47").Value = '1.957.63'
$ws.Range("E47").Value = '  -9.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0221'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0875'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -10.06%  '
